$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a plain text value to a cell without letting Excel's COM
# layer auto-coerce numeric-looking strings (e.g. "27.28") into numbers.
# Temporarily force the cell to Text format, assign, then restore the
# cell's style back to the default "Normal" so no stray number-format
# styling is left behind on the cell.
function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue "D2" "63.521.19"
Set-TextValue "E2" "  +2.66%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.477.00"
Set-TextValue "E3" "  +2.57%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  +0.24%  "

# Row 5 - BNB
Set-TextValue "D5" "573.88"
Set-TextValue "E5" "  +2.07%  "

# Row 6 - Solana
Set-TextValue "D6" "149.88"
Set-TextValue "E6" "  +5.30%  "

# Row 7 - USDC
Set-TextValue "E7" "  -0.13%  "

# Row 8 - XRP
Set-TextValue "E8" "  +2.23%  "

# Row 10 - TRON
Set-TextValue "E10" "  +0.47%  "

# Row 11 - Cardano
Set-TextValue "E11" "  +4.63%  "

# Row 12 - Toncoin
Set-TextValue "E12" "  +2.65%  "

# Row 13 - Avalanche
Set-TextValue "D13" "27.28"
Set-TextValue "E13" "  +6.19%  "

# Row 14 - ShibaInu
Set-TextValue "E14" "  +7.43%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "E15" "  +2.17%  "

# Row 16 - WrappedBTC
Set-TextValue "D16" "63.361.44"
Set-TextValue "E16" "  +2.59%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.482.21"
Set-TextValue "E17" "  +2.16%  "

# Row 18 - Chainlink
Set-TextValue "D18" "11.60"
Set-TextValue "E18" "  +2.98%  "

# Row 19 - Uniswap
Set-TextValue "D19" "7.26"
Set-TextValue "E19" "  +7.10%  "

# Row 20 - Polkadot
Set-TextValue "E20" "  +3.49%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "328.92"
Set-TextValue "E21" "  +1.83%  "

# Row 23 - SuiNetwork
Set-TextValue "E23" "  +11.18%  "

# Row 24 - Litecoin
Set-TextValue "D24" "67.60"
Set-TextValue "E24" "  +1.74%  "

# Row 25 - Bittensor
Set-TextValue "D25" "638.44"
Set-TextValue "E25" "  +15.34%  "

# Row 26 - Aptos
Set-TextValue "E26" "  +1.36%  "

# Row 27 - PEPE
Set-TextValue "E27" "  +13.24%  "

# Row 28 - WrappedeETH
Set-TextValue "D28" "2.604.41"
Set-TextValue "E28" "  +2.66%  "

# Row 29 - Fetch.AI
Set-TextValue "E29" "  +10.07%  "

# Row 30 - InternetComputer(DFINITY)
Set-TextValue "D30" "8.50"
Set-TextValue "E30" "  +4.34%  "

# Row 31 - Binance-PegBSC-USD
Set-TextValue "D31" "0.992"
Set-TextValue "E31" "  -0.94%  "

# Row 32 - Kaspa
Set-TextValue "E32" "  -1.70%  "

# Row 33 - PancakeSwap
Set-TextValue "E33" "  +3.20%  "

# Row 34 - NEARProtocol
Set-TextValue "D34" "5.24"
Set-TextValue "E34" "  +10.98%  "

# Row 35 - ImmutableX
Set-TextValue "E35" "  +5.17%  "

# Row 36 - FirstDigitalUSD
Set-TextValue "D36" "0.998"
Set-TextValue "E36" "  -0.16%  "

# Row 37 - PolygonEcosystemToken
Set-TextValue "E37" "  +2.40%  "

# Row 38 - RenderToken
Set-TextValue "E38" "  +2.00%  "

# Row 39 - EthereumClassic
Set-TextValue "D39" "18.96"
Set-TextValue "E39" "  +2.54%  "

# Row 40 - Stacks
Set-TextValue "E40" "  +2.68%  "

# Row 41 - Monero
Set-TextValue "D41" "146.82"
Set-TextValue "E41" "  -4.52%  "

# Row 42 - dogwifhat
Set-TextValue "E42" "  +18.47%  "

# Row 43 - USDe
Set-TextValue "E43" "  +0.81%  "

# Row 44 - Aave
Set-TextValue "D44" "151.62"
Set-TextValue "E44" "  +3.43%  "

# Row 45 - Filecoin
Set-TextValue "D45" "3.79"
Set-TextValue "E45" "  +4.90%  "

# Row 46 - Hedera
Set-TextValue "D46" "0.0554"
Set-TextValue "E46" "  +5.64%  "

# Row 47 - InjectiveProtocol
Set-TextValue "D47" "21.15"
Set-TextValue "E47" "  +7.41%  "

# Row 48 - Mantle
Set-TextValue "E48" "  +3.56%  "

# Row 49 - VeChain
Set-TextValue "D49" "0.0241"
Set-TextValue "E49" "  +6.63%  "

# Row 50 - Stellar
Set-TextValue "D50" "0.0928"
Set-TextValue "E50" "  +1.19%  "

# Row 51 - WhiteBITCoin -> ONDO (coin replaced entirely)
Set-TextValue "B51" "ONDO"
Set-TextValue "C51" "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
Set-TextValue "D51" "0.743"
Set-TextValue "E51" "  +5.90%  "
